# feat: add 2022-Q1 data
#
# Plan:
#  1. Duplicate the existing "总计" (totals) sheet to the end of the workbook
#     BEFORE touching it - the duplicate will become the new "总计" sheet
#     (with an extra row for 2022-Q1), while the original sheet (still in
#     position 6, right after "2021-Q4") gets renamed to "2022-Q1" and
#     repopulated with the new quarter's per-fund holdings.
#  2. Rename + repopulate the original "总计" sheet -> "2022-Q1".
#  3. Rename the duplicated sheet -> "总计", insert a new row at the top for
#     "2022-Q1" and fill in the numbers, renumbering the helper index column.
#  4. Restore the original active sheet/selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 0: duplicate "总计" to the end of the workbook (keeps exact formatting)
# ---------------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))

# Sheet objects, resolved fresh by position (names collide until we rename).
$newQuarterSheet = $wb.Worksheets.Item(6)   # was "总计", becomes "2022-Q1"
$newTotalSheet   = $wb.Worksheets.Item(7)   # the copy, becomes "总计"

# ---------------------------------------------------------------------------
# Step 1: "2022-Q1" fund-holdings sheet (replaces the old "总计" sheet data)
# ---------------------------------------------------------------------------
$newQuarterSheet.Name = "2022-Q1"
$newQuarterSheet.Cells.Clear()

# Use "2021-Q4" as the formatting template - same header/style layout.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("A1:H4").Copy()
$newQuarterSheet.Range("A1").PasteSpecial(-4163)   # xlPasteValues
$template.Range("A1:H4").Copy()
$newQuarterSheet.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$newQuarterSheet.Range("A1").Clear()               # A1 itself stays empty
$newQuarterSheet.Range("A3:H4").Delete()           # only one data row this quarter

# Fill in the real 2022-Q1 values (row 2). Values that look numeric
# ("160613", "3.05", ...) must stay plain text (matches the source data,
# which stores every column but H as text) - format the cell as Text first
# so the COM layer doesn't silently coerce the string into a Double, then
# drop the throwaway NumberFormat style so the cell is back to the default,
# un-styled look (matching its neighbours).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $newQuarterSheet.Range("B2") "160613"
Set-TextValue $newQuarterSheet.Range("C2") "鹏华盛世创新混合(LOF)"
Set-TextValue $newQuarterSheet.Range("D2") "3.05"
Set-TextValue $newQuarterSheet.Range("E2") "93.82"
Set-TextValue $newQuarterSheet.Range("F2") "3.19"
Set-TextValue $newQuarterSheet.Range("G2") "0.0973"
$newQuarterSheet.Range("H2").Value = 10

# ---------------------------------------------------------------------------
# Step 2: "总计" sheet - insert the 2022-Q1 summary row at the top
# ---------------------------------------------------------------------------
$newTotalSheet.Name = "总计"
$newTotalSheet.Rows("2").Insert()

# Row insert copies the style of the row above into the new blank row
# (bordered "s=2" look) - but only column A should keep that style, so
# reset B:D back to the default look before filling them in.
$newTotalSheet.Range("B2:D2").ClearFormats()

$newTotalSheet.Range("A6").Copy()
$newTotalSheet.Range("A2").PasteSpecial(-4122)     # xlPasteFormats (style 2)
$newTotalSheet.Range("A2").Value = 0

$newTotalSheet.Range("B2").Value = "2022-Q1"
$newTotalSheet.Range("C2").Value = 1
$newTotalSheet.Range("D2").Value = 0.1

# Renumber the helper index column (A) for the rows that shifted down.
$newTotalSheet.Range("A3").Value = 1
$newTotalSheet.Range("A4").Value = 2
$newTotalSheet.Range("A5").Value = 3
$newTotalSheet.Range("A6").Value = 4
$newTotalSheet.Range("A7").Value = 5

# ---------------------------------------------------------------------------
# Step 3: restore original selection (first sheet active)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
